$d = $word.ActiveDocument

# The 2nd paragraph in the original document is empty, styled "Heading 1",
# and carries the hidden _GoBack bookmark at its (only) position. We insert
# 7 new empty paragraphs immediately before it so that it remains the final
# (9th) paragraph of the document -- keeping the bookmark exactly where the
# target diff wants it (inside the very last paragraph).
$lastP = $d.Paragraphs.Item(2)
for ($i = 0; $i -lt 7; $i++) {
    $lastP.Range.InsertParagraphBefore() | Out-Null
}

# Appends text (optionally split into runs, some italic) to the end of a
# paragraph, i.e. right before its trailing paragraph mark. Using an
# explicit "insertion point" range anchored at (End-1) avoids accidentally
# spilling text into the *next* paragraph, which is what Range.End actually
# points at (the position right after the paragraph mark == next para's
# start).
function Add-Segments($para, $segments) {
    foreach ($seg in $segments) {
        $text = $seg[0]
        $italic = $seg[1]
        $insPos = $para.Range.End - 1
        $ip = $d.Range($insPos, $insPos)
        $ip.InsertAfter($text)
        if ($italic) {
            $run = $d.Range($insPos, $insPos + $text.Length)
            $run.Font.Italic = 1
        }
    }
}

function New-Heading($para, [string]$text) {
    $para.Style = "Heading 1"
    $para.Range.Text = $text
}

function New-Body($para, [bool]$justify, $segments) {
    $para.Style = "Normal"
    if ($justify) {
        $para.Alignment = 3
    }
    $para.Range.Text = ""
    Add-Segments $para $segments
}

# --- Paragraph 2: "Puntos de interés" heading ---
New-Heading $d.Paragraphs.Item(2) "Puntos de interés"

# --- Paragraph 3 ---
New-Body $d.Paragraphs.Item(3) $true @(
    ,("Los puntos de interés decidimos manejarlos con una herencia con una clase padre ", $false)
    ,("POI", $true)
    ,(" que define los métodos genéricos de un punto de interés, ya que los mismos se redefinen en cada clase hija que representa cada especificación de un punto de interés.", $false)
)

# --- Paragraph 4 ---
New-Body $d.Paragraphs.Item(4) $true @(
    ,("Cada una de las especificaciones redefine los métodos en los cuales el método no es genérico y añaden los campos que necesitan.", $false)
)

# --- Paragraph 5: "Repositorio Local" heading ---
New-Heading $d.Paragraphs.Item(5) "Repositorio Local"

# --- Paragraph 6 ---
New-Body $d.Paragraphs.Item(6) $false @(
    ,("El repositorio local lo creamos como una clase hija de ", $false)
    ,("CollectionBasedRepo", $true)
    ,(" por lo cual tuvimos que hacer que ", $false)
    ,("POI", $true)
    ,(" herede de la clase ", $false)
    ,("Entity", $true)
    ,(" para poder ingresarlos al repositorio (hecho que nos forzó a hacer concreta la clase ", $false)
    ,("POI", $true)
    ,(" y no abstracta como antes). En ", $false)
    ,("POI", $true)
    ,(" se redefinieron algunos métodos necesarios para validación heredados de Entity.", $false)
)

# --- Paragraph 7: "Repositorio Externo" heading ---
New-Heading $d.Paragraphs.Item(7) "Repositorio Externo"

# --- Paragraph 8 ---
New-Body $d.Paragraphs.Item(8) $false @(
    ,("Para el manejo del servicio externo tuvimos que crear una clase ", $false)
    ,("AdaptadorServicioExterno", $true)
    ,(" que desencripte los resultados del servicio externo representado por la interfaz ", $false)
    ,("interfazConsultaBancaria", $true)
    ,(" (para los tests se usó un mock de la clase). Además, todos los orígenes de datos (Adaptador y Repositorio Local) implementan la interfaz ", $false)
    ,("OrigenDeDatos", $true)
    ,(" que define como debe ser una consulta genérica.", $false)
)

# --- Paragraph 9: original (2nd) paragraph; keeps the _GoBack bookmark ---
New-Body $d.Paragraphs.Item(9) $false @(
    ,("El funcionamiento del mismo primero fue mediante una clase Stub para el servicio externo pero el mismo no era polimórfico para distintos orígenes de datos, por lo que se decidió que el repositorio local tenga un conjunto de servicios externos, entonces si un POI no era encontrado en el repositorio local se iba a buscar a el/los repositorios externos.", $false)
)
